$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the DATA_HORA_ANALISE timestamp (stored as text)
$ws.Range("A2").Value = "2025-05-25 18:32:19"

# Update the metric values in row 2
$ws.Range("B2").Value = 16229
$ws.Range("C2").Value = 11751
$ws.Range("D2").Value = 72.40741881816501
$ws.Range("E2").Value = 2241
$ws.Range("F2").Value = 13.8086142091318
$ws.Range("G2").Value = 3105
$ws.Range("H2").Value = 19.13241727771274
$ws.Range("I2").Value = 9558
$ws.Range("J2").Value = 58.8945714461766
$ws.Range("K2").Value = 3074903.87
$ws.Range("L2").Value = 3566
$ws.Range("M2").Value = 21.97301127611066
$ws.Range("N2").Value = 1158493.43
$ws.Range("O2").Value = 4689
$ws.Range("P2").Value = 28.89272290344445
$ws.Range("Q2").Value = 494295.3
$ws.Range("R2").Value = 3569
$ws.Range("S2").Value = 21.99149670343213
$ws.Range("T2").Value = 3435
$ws.Range("U2").Value = 21.16581428307351
$ws.Range("V2").Value = 2427830.57
$ws.Range("W2").Value = 2299
$ws.Range("X2").Value = 14.16599913734672
$ws.Range("Y2").Value = 1434
$ws.Range("Z2").Value = 8.836034259658636
$ws.Range("AA2").Value = 152778
$ws.Range("AB2").Value = 803
$ws.Range("AC2").Value = 4.947932713044549
$ws.Range("AD2").Value = 16241
$ws.Range("AE2").Value = 10611
$ws.Range("AF2").Value = 65.33464688134967
$ws.Range("AG2").Value = 5630
$ws.Range("AH2").Value = 34.66535311865033
$ws.Range("AI2").Value = 479
$ws.Range("AJ2").Value = 866
$ws.Range("AK2").Value = 1354
$ws.Range("AL2").Value = 17.74731381993331
$ws.Range("AM2").Value = 32.08595776213413
$ws.Range("AN2").Value = 50.16672841793257
$ws.Range("AO2").Value = 1360191
$ws.Range("AP2").Value = 255433.51
$ws.Range("AQ2").Value = 85168.91
$ws.Range("AR2").Value = 79.97391005898884
$ws.Range("AS2").Value = 15.01849119336316
$ws.Range("AT2").Value = 5.007598747648025
$ws.Range("AU2").Value = 49.75165717486696
$ws.Range("AV2").Value = 239.0586868099942
$ws.Range("AW2").Value = 535.3431192660551

$wb.Save()
